$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A150:A171").NumberFormat = "@"

$ws.Range("A150").Value = "05-08-2021"
$ws.Range("B150").Value = 13208
$ws.Range("C150").Value = 21450
$ws.Range("D150").Value = -8242
$ws.Range("A151").Value = "06-08-2021"
$ws.Range("B151").Value = 13265
$ws.Range("C151").Value = 21622
$ws.Range("D151").Value = -8357
$ws.Range("A152").Value = "09-08-2021"
$ws.Range("B152").Value = 13793
$ws.Range("C152").Value = 21919
$ws.Range("D152").Value = -8126
$ws.Range("A153").Value = "10-08-2021"
$ws.Range("B153").Value = 13315
$ws.Range("C153").Value = 21444
$ws.Range("D153").Value = -8129
$ws.Range("A154").Value = "11-08-2021"
$ws.Range("B154").Value = 13570
$ws.Range("C154").Value = 21554
$ws.Range("D154").Value = -7984
$ws.Range("A155").Value = "12-08-2021"
$ws.Range("B155").Value = 13517
$ws.Range("C155").Value = 21531
$ws.Range("D155").Value = -8014
$ws.Range("A156").Value = "13-08-2021"
$ws.Range("B156").Value = 13705
$ws.Range("C156").Value = 21318
$ws.Range("D156").Value = -7614
$ws.Range("A157").Value = "16-08-2021"
$ws.Range("B157").Value = 13485
$ws.Range("C157").Value = 21214
$ws.Range("D157").Value = -7729
$ws.Range("A158").Value = "17-08-2021"
$ws.Range("B158").Value = 13325
$ws.Range("C158").Value = 21040
$ws.Range("D158").Value = -7715
$ws.Range("A159").Value = "18-08-2021"
$ws.Range("B159").Value = 13354
$ws.Range("C159").Value = 20966
$ws.Range("D159").Value = -7612
$ws.Range("A160").Value = "19-08-2021"
$ws.Range("B160").Value = 13613
$ws.Range("C160").Value = 21030
$ws.Range("D160").Value = -7417
$ws.Range("A161").Value = "20-08-2021"
$ws.Range("B161").Value = 13460
$ws.Range("C161").Value = 20886
$ws.Range("D161").Value = -7426
$ws.Range("A162").Value = "23-08-2021"
$ws.Range("B162").Value = 13440
$ws.Range("C162").Value = 20987
$ws.Range("D162").Value = -7547
$ws.Range("A163").Value = "24-08-2021"
$ws.Range("B163").Value = 13567
$ws.Range("C163").Value = 20925
$ws.Range("D163").Value = -7358
$ws.Range("A164").Value = "25-08-2021"
$ws.Range("B164").Value = 13423
$ws.Range("C164").Value = 20814
$ws.Range("D164").Value = -7391
$ws.Range("A165").Value = "26-08-2021"
$ws.Range("B165").Value = 13330
$ws.Range("C165").Value = 20542
$ws.Range("D165").Value = -7212
$ws.Range("A166").Value = "27-08-2021"
$ws.Range("B166").Value = 13178
$ws.Range("C166").Value = 20263
$ws.Range("D166").Value = -7085
$ws.Range("A167").Value = "30-08-2021"
$ws.Range("B167").Value = 13538
$ws.Range("C167").Value = 20305
$ws.Range("D167").Value = -6767
$ws.Range("A168").Value = "31-08-2021"
$ws.Range("B168").Value = 13679
$ws.Range("C168").Value = 19961
$ws.Range("D168").Value = -6283
$ws.Range("A169").Value = "01-09-2021"
$ws.Range("B169").Value = 14126
$ws.Range("C169").Value = 19616
$ws.Range("D169").Value = -5491
$ws.Range("A170").Value = "02-09-2021"
$ws.Range("B170").Value = 14511
$ws.Range("C170").Value = 19618
$ws.Range("D170").Value = -5107
$ws.Range("A171").Value = "03-09-2021"
$ws.Range("B171").Value = 15083
$ws.Range("C171").Value = 19716
$ws.Range("D171").Value = -4633

$ws.Range("A150:A171").ClearFormats()
